# Rename parameter symbol names in column D from snake_case to camelCase.
# Also fixes v_static/duration_static -> vStop/durationStop (renamed, not just re-cased).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$renames = @{
    "D3"  = "ifRecalib"
    "D4"  = "calibSeconds(s)"
    "D5"  = "cellLen(m)"
    "D6"  = "qMerge(v/h)"
    "D7"  = "laneWidth(m)"
    "D8"  = "emgcWidth(m)"
    "D9"  = "qCalDuration(s)"
    "D10" = "calInterval(s)"
    "D11" = "maxCompleteFrames"
    "D12" = "smoothAlpha"
    "D13" = "eventTypes"
    "D14" = "tTolerance(s)"
    "D15" = "qStandard(v/h)"
    "D18" = "vLateral(m/s)"
    "D19" = "spillWarnFrequecy(s)"
    "D20" = "vStop(m/s)"
    "D21" = "durationStop(s)"
    "D22" = "vLow(m/s)"
    "D23" = "durationLow(s)"
    "D24" = "vHigh(m/s)"
    "D25" = "durationHigh(s)"
    "D26" = "aIntense(m/s^2)"
    "D27" = "durationIntense(s)"
    "D28" = "dTouch(m)"
    "D29" = "tSupervise(s)"
    "D30" = "densityCrowd(pcu/km/ln)"
    "D31" = "vCrowd(m/s)"
    "D32" = "durationOccupation(s)"
}

foreach ($cellRef in $renames.Keys) {
    $ws.Range($cellRef).Value = $renames[$cellRef]
}

# Match final cell selection left behind in the saved file.
$ws.Range("D33").Select()
